# Replace space in group column names on the "group_cooking" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("group_cooking")

$ws.Range("A1").Value = "Cooking_Equipment"
$ws.Range("B1").Value = "Years_Owned"

# Activate this sheet and move the selection to B2, matching the author's
# final cursor position after editing the header cells.
$ws.Activate()
$ws.Range("B2").Select()
